$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.347.94'
$ws.Range("E2").Value = '  +4.02%  '
$ws.Range("D3").Value = '3.486.29'
$ws.Range("E3").Value = '  +3.48%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '584.28'
$ws.Range("E5").Value = '  +2.20%  '
$ws.Range("D6").Value = '147.57'
$ws.Range("E6").Value = '  +6.77%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = '0.477'
$ws.Range("E8").Value = '  +1.18%  '
$ws.Range("D9").Value = '7.69'
$ws.Range("E9").Value = '  -0.08%  '
$ws.Range("E10").Value = '  +4.36%  '
$ws.Range("D11").Value = '0.397'
$ws.Range("E11").Value = '  +3.83%  '
$ws.Range("D12").Value = '4.082.79'
$ws.Range("E12").Value = '  +3.54%  '
$ws.Range("D13").Value = '29.68'
$ws.Range("E13").Value = '  +5.70%  '
$ws.Range("E14").Value = '  -0.42%  '
$ws.Range("D15").Value = '3.486.23'
$ws.Range("E15").Value = '  +3.73%  '
$ws.Range("E16").Value = '  +3.38%  '
$ws.Range("D17").Value = '63.361.97'
$ws.Range("E17").Value = '  +3.91%  '
$ws.Range("D18").Value = '6.28'
$ws.Range("E18").Value = '  +3.07%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.40'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +6.60%  '
$ws.Range("D20").Value = '9.38'
$ws.Range("E20").Value = '  +5.38%  '
$ws.Range("D21").Value = '391.31'
$ws.Range("E21").Value = '  +1.75%  '
$ws.Range("E22").Value = '  +2.40%  '
$ws.Range("D23").Value = '75.28'
$ws.Range("E23").Value = '  +0.38%  '
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("E25").Value = '  +7.96%  '
$ws.Range("D26").Value = '3.628.29'
$ws.Range("E26").Value = '  +3.56%  '
$ws.Range("E27").Value = '  -2.38%  '
$ws.Range("E28").Value = '  +9.99%  '
$ws.Range("E29").Value = '  +0.15%  '
$ws.Range("D30").Value = '8.28'
$ws.Range("E30").Value = '  +4.47%  '
$ws.Range("E31").Value = '  +1.87%  '
$ws.Range("E32").Value = '  +7.51%  '
$ws.Range("D34").Value = '23.78'
$ws.Range("E34").Value = '  +3.30%  '
$ws.Range("D35").Value = '32.59'
$ws.Range("E35").Value = '  +28.18%  '
$ws.Range("D36").Value = '5.33'
$ws.Range("E36").Value = '  +8.49%  '
$ws.Range("D37").Value = '7.13'
$ws.Range("E37").Value = '  +4.62%  '
$ws.Range("D38").Value = '171.54'
$ws.Range("E38").Value = '  +2.61%  '
$ws.Range("E39").Value = '  +9.01%  '
$ws.Range("D40").Value = '3.522.72'
$ws.Range("E40").Value = '  +3.44%  '
$ws.Range("D41").Value = '0.0767'
$ws.Range("E41").Value = '  +1.59%  '
$ws.Range("D42").Value = '0.807'
$ws.Range("E42").Value = '  +4.72%  '
$ws.Range("E43").Value = '  +3.64%  '
$ws.Range("D44").Value = '42.45'
$ws.Range("E44").Value = '  +0.43%  '
$ws.Range("E45").Value = '  +6.68%  '
$ws.Range("E46").Value = '  +9.46%  '
$ws.Range("D47").Value = '2.629.77'
$ws.Range("E47").Value = '  +7.76%  '
$ws.Range("D48").Value = '23.66'
$ws.Range("E48").Value = '  +7.35%  '
$ws.Range("E49").Value = '  +16.20%  '
$ws.Range("D50").Value = '6.75'
$ws.Range("E50").Value = '  +1.68%  '
$ws.Range("E51").Value = '  +5.30%  '
